# Applies the "Added Final Report & PowerPoint" edit:
#  1. "Group #4/5 Project #2..." -> "Group #5 Project #2..."
#  2. Merge the "Steam "/"G"/"aming"/" Analytics" runs into a single
#     "Steam Gaming Analytics" run (same visible text, tidied markup).
#  3. The 7th blank centered paragraph following the "Steam Gaming
#     Analytics" heading becomes right-aligned instead of centered.

$d = $word.ActiveDocument

# 1) "4/5" -> "5" (unique occurrence in the document)
$d.Content.Find.Execute("4/5", $false, $false, $false, $false, $false, $true, 1, $false, "5", 2) | Out-Null

# 2) Consolidate the "Steam Gaming Analytics" heading into one run.
$d.Content.Find.Execute("Steam Gaming Analytics", $false, $false, $false, $false, $false, $true, 1, $false, "Steam Gaming Analytics", 2) | Out-Null

# 3) Find the "Steam Gaming Analytics" heading paragraph, then walk
#    forward 7 paragraphs (the blank, centered ones) and switch that
#    paragraph's alignment from center to right.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "Steam Gaming Analytics") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $p = $target
    for ($i = 0; $i -lt 7; $i++) {
        $p = $p.Next()
    }
    $p.Alignment = 2
}
